# Update the "Estado de Cuenta" worksheet: rotate the employee data rows
# (16-18) so that the newest arrear record (JHON JAIRO ESTEVEZ ROMERO,
# 19601056) appears first, followed by the previously-first two records.
# Row 19 (second period for the same worker) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: now holds what used to be row 18's data
$ws.Range("C16").Value = "19601056"
$ws.Range("D16").Value = "JHON JAIRO ESTEVEZ ROMERO"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Row 17: now holds what used to be row 16's data
$ws.Range("C17").Value = "1143389778"
$ws.Range("D17").Value = "DAVID ABRAHAM DOVALE LUNA"
$ws.Range("F17").Value = 12096
$ws.Range("G17").Value = 1512000

# Row 18: now holds what used to be row 17's data
$ws.Range("C18").Value = "76330112"
$ws.Range("D18").Value = "HECTOR FABIO ARCOS SANCHEZ"
$ws.Range("F18").Value = 11388
$ws.Range("G18").Value = 1423500

$wb.Save()
